# Applies the scheduled-runner profit recalculation update to the
# Sheets workbook (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR), row-by-row,
# matching the authoritative OOXML diff cell-for-cell.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 10666.667
$ws.Range("J3").Value = 10666.667
$ws.Range("L3").Value = 10666.667
$ws.Range("N3").Value = -10894.667
$ws.Range("H10").Value = 4399.933
$ws.Range("J10").Value = 4999.9165
$ws.Range("L10").Value = 4999.9165
$ws.Range("N10").Value = -5585.9165
$ws.Range("H102").Value = 10666.667
$ws.Range("J102").Value = 10666.667
$ws.Range("L102").Value = 10666.667
$ws.Range("N102").Value = -17156.667

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 117
$ws.Range("J5").Value = 134
$ws.Range("L5").Value = 134
$ws.Range("N5").Value = -358
$ws.Range("H11").Value = 5006000
$ws.Range("I11").Value = 10000000
$ws.Range("J11").Value = 12000
$ws.Range("K11").Value = 10000000
$ws.Range("L11").Value = 12000
$ws.Range("M11").Value = -9999856
$ws.Range("N11").Value = -12288
$ws.Range("H24").Value = 36500
$ws.Range("J24").Value = 36500
$ws.Range("L24").Value = 36500
$ws.Range("N24").Value = -37248
$ws.Range("H32").Value = 3363.4285
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H61").Value = 4495.8335
$ws.Range("I61").Value = 3417.4285
$ws.Range("K61").Value = 3417.4285
$ws.Range("M61").Value = -3205.4285
$ws.Range("H100").Value = 36500
$ws.Range("J100").Value = 36500
$ws.Range("L100").Value = 36500
$ws.Range("N100").Value = -38664
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H136").Value = 4495.8335
$ws.Range("I136").Value = 3417.4285
$ws.Range("K136").Value = 10252.2855
$ws.Range("M136").Value = -7702.2855

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 117
$ws.Range("J4").Value = 134
$ws.Range("L4").Value = 134
$ws.Range("N4").Value = -364
$ws.Range("H107").Value = 1545.0577
$ws.Range("I107").Value = 1528.0294
$ws.Range("J107").Value = 1577.2222
$ws.Range("K107").Value = 1528.0294
$ws.Range("L107").Value = 1577.2222
$ws.Range("M107").Value = 391.9706000000001
$ws.Range("N107").Value = -5417.2222

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1228.1818
$ws.Range("I10").Value = 1228.1818
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 1228.1818
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -1089.1818
$ws.Range("N10").ClearContents()
$ws.Range("H31").Value = 1285.8667
$ws.Range("I31").Value = 1086.5
$ws.Range("J31").Value = 2083.3333
$ws.Range("K31").Value = 1086.5
$ws.Range("L31").Value = 2083.3333
$ws.Range("M31").Value = -791.5
$ws.Range("N31").Value = -2673.3333
$ws.Range("H34").Value = 1285.8667
$ws.Range("I34").Value = 1086.5
$ws.Range("J34").Value = 2083.3333
$ws.Range("K34").Value = 1086.5
$ws.Range("L34").Value = 2083.3333
$ws.Range("M34").Value = -884.5
$ws.Range("N34").Value = -2487.3333
$ws.Range("H35").Value = 1220.75
$ws.Range("I35").Value = 1240.8182
$ws.Range("J35").Value = 1000
$ws.Range("K35").Value = 1240.8182
$ws.Range("L35").Value = 1000
$ws.Range("M35").Value = -946.8181999999999
$ws.Range("N35").Value = -1588
$ws.Range("H86").Value = 2121.05
$ws.Range("I86").Value = 1676
$ws.Range("J86").Value = 2485.182
$ws.Range("K86").Value = 1676
$ws.Range("L86").Value = 2485.182
$ws.Range("M86").Value = -553
$ws.Range("N86").Value = -4731.182
$ws.Range("H89").Value = 2121.05
$ws.Range("I89").Value = 1676
$ws.Range("J89").Value = 2485.182
$ws.Range("K89").Value = 8380
$ws.Range("L89").Value = 12425.91
$ws.Range("M89").Value = -2764
$ws.Range("N89").Value = -23657.91

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 370.85715
$ws.Range("I17").Value = 299.5
$ws.Range("J17").Value = 399.4
$ws.Range("K17").Value = 898.5
$ws.Range("L17").Value = 1198.2
$ws.Range("M17").Value = -729.5
$ws.Range("N17").Value = -1536.2
$ws.Range("H98").Value = 750.6667
$ws.Range("I98").Value = 720
$ws.Range("J98").Value = 904
$ws.Range("K98").Value = 2160
$ws.Range("L98").Value = 2712
$ws.Range("M98").Value = -662
$ws.Range("N98").Value = -5708
$ws.Range("H139").Value = 2069.1667
$ws.Range("I139").Value = 1603.75
$ws.Range("K139").Value = 4811.25
$ws.Range("M139").Value = 328.75
$ws.Range("H140").Value = 1097.375
$ws.Range("I140").Value = 970.5333000000001
$ws.Range("J140").Value = 3000
$ws.Range("K140").Value = 2911.5999
$ws.Range("L140").Value = 9000
$ws.Range("M140").Value = 2268.4001
$ws.Range("N140").Value = -19360

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 431177.78
$ws.Range("J21").Value = 2001329.6
$ws.Range("L21").Value = 2001329.6
$ws.Range("N21").Value = -2001675.6
$ws.Range("H30").Value = 431177.78
$ws.Range("J30").Value = 2001329.6
$ws.Range("L30").Value = 2001329.6
$ws.Range("N30").Value = -2001539.6
$ws.Range("H57").Value = 14598
$ws.Range("J57").Value = 15108.889
$ws.Range("L57").Value = 15108.889
$ws.Range("N57").Value = -16748.889
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H132").Value = 3207.1304
$ws.Range("I132").Value = 2492.1538
$ws.Range("J132").Value = 4136.6
$ws.Range("K132").Value = 7476.4614
$ws.Range("L132").Value = 12409.8
$ws.Range("M132").Value = -4946.4614
$ws.Range("N132").Value = -17469.8

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H101").Value = 17590.25
$ws.Range("J101").Value = 17590.25
$ws.Range("L101").Value = 17590.25
$ws.Range("N101").Value = -24080.25
$ws.Range("H104").Value = 15375
$ws.Range("J104").Value = 15375
$ws.Range("L104").Value = 15375
$ws.Range("N104").Value = -22363

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H104").Value = 20481.285
$ws.Range("J104").Value = 20481.285
$ws.Range("L104").Value = 20481.285
$ws.Range("N104").Value = -27469.285
$ws.Range("H105").Value = 35000
$ws.Range("J105").Value = 35000
$ws.Range("L105").Value = 35000
$ws.Range("N105").Value = -41988
